$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# REPORTGEN-709: the "Application Name:" / "Snapshot date:" label cells in
# B3/B4 used to be merged with C3/C4, with the RepGen placeholder tag
# living one column further right (D3/D4). Unmerge the label cells and
# slide each RepGen tag one column to the left (D3 -> C3, D4 -> C4),
# leaving D3/D4 blank but formatted exactly as they were before.

$ws.Range("B3:C3").UnMerge()
$ws.Range("D3").Copy($ws.Range("C3"))
$ws.Range("D3").Value2 = $null

$ws.Range("B4:C4").UnMerge()
$ws.Range("D4").Copy($ws.Range("C4"))
$ws.Range("D4").Value2 = $null

$excel.CutCopyMode = $false

$ws.Activate()
$ws.Range("D3:D4").Select()

Write-Host "REPORTGEN-709: update full detailed excel reports - Summary sheet layout updated"
